# Penalty Reward System tweak — shift the Week_Start_Date by one week and
# bump up the MyForecast numbers on the "Forecast Comparison" sheet, then
# refresh the dependent stats on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Helper: write a literal text value into a cell without leaving it tagged
# as a date/number-formatted cell (Excel auto-detects "2025-01-12"-style
# strings as dates otherwise). A leading apostrophe forces text entry, and
# resetting the style afterwards keeps the cell's formatting identical to
# its original (unstyled) state.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# --- Forecast Comparison sheet -------------------------------------------
# Row -> (new Week_Start_Date, new MyForecast)
$forecastUpdates = @(
    @{ Row = 2;  Date = "2025-01-12"; Forecast = 113 },
    @{ Row = 3;  Date = "2025-01-19"; Forecast = 117 },
    @{ Row = 4;  Date = "2025-01-26"; Forecast = 118 },
    @{ Row = 5;  Date = "2025-02-02"; Forecast = 117 },
    @{ Row = 6;  Date = "2025-02-09"; Forecast = 120 },
    @{ Row = 7;  Date = "2025-02-16"; Forecast = 127 },
    @{ Row = 8;  Date = "2025-02-23"; Forecast = 131 },
    @{ Row = 9;  Date = "2025-03-02"; Forecast = 130 },
    @{ Row = 10; Date = "2025-03-09"; Forecast = 127 },
    @{ Row = 11; Date = "2025-03-16"; Forecast = 127 },
    @{ Row = 12; Date = "2025-03-23"; Forecast = 134 },
    @{ Row = 13; Date = "2025-03-30"; Forecast = 142 },
    @{ Row = 14; Date = "2025-04-06"; Forecast = 148 },
    @{ Row = 15; Date = "2025-04-13"; Forecast = 150 },
    @{ Row = 16; Date = "2025-04-20"; Forecast = 153 },
    @{ Row = 17; Date = "2025-04-27"; Forecast = 154 }
)

foreach ($u in $forecastUpdates) {
    Set-TextValue $wsForecast.Cells.Item($u.Row, 2) $u.Date
    $wsForecast.Cells.Item($u.Row, 4).Value = $u.Forecast
}

# --- Summary sheet ---------------------------------------------------------
Set-TextValue $wsSummary.Cells.Item(2, 2)  "2023-11-19 to 2025-01-05"
Set-TextValue $wsSummary.Cells.Item(4, 2)  "180"
Set-TextValue $wsSummary.Cells.Item(5, 2)  "78"
Set-TextValue $wsSummary.Cells.Item(7, 2)  "46"
Set-TextValue $wsSummary.Cells.Item(8, 2)  "4622 units"
Set-TextValue $wsSummary.Cells.Item(9, 2)  "2108"
Set-TextValue $wsSummary.Cells.Item(10, 2) "973"
Set-TextValue $wsSummary.Cells.Item(11, 2) "465"
Set-TextValue $wsSummary.Cells.Item(12, 2) "154"
Set-TextValue $wsSummary.Cells.Item(13, 2) "2025-04-27"
Set-TextValue $wsSummary.Cells.Item(14, 2) "113"
Set-TextValue $wsSummary.Cells.Item(15, 2) "2025-01-12"
